$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.510.88'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.02%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.871.13'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.13%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '466.96'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +9.91%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '148.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +13.15%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.633'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +3.18%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.754'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.17%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.156'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.36%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000314'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -8.88%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.88'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.01%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.45'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.73%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.480.75'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.34%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.80'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -7.08%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.890.88'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.68%  '

$ws.Range('E17').Value = '  -0.26%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '20.10'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.63%  '

$ws.Range('E19').Value = '  +7.46%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '67.566.72'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.96%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '433.48'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.86%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '14.84'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.40%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.29'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +8.46%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.70'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.98%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +9.86%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.30'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +13.22%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '37.67'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.29%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.20'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.66%  '

$ws.Range('E29').Value = '  +4.77%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '743.64'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.45%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.136'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +10.71%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '13.76'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.46%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.76'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.71%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '43.17'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +10.40%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.164'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +7.37%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '57.73'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.90%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.12%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '5.53'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.64%  '

$ws.Range('E39').Value = '  +3.76%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.351'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +11.90%  '

$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.92'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.24%  '

$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.63'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +15.69%  '

$ws.Range('E43').Value = '  +5.46%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0₃0682'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -10.11%  '

$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.44'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.37%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.25'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.44%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.77'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +7.81%  '

$ws.Range('E49').Value = '  +3.64%  '

$ws.Range('E50').Value = '  +3.19%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '144.18'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.92%  '
